# Refresh crypto Price (D) and Volume(1h) (E) columns with updated market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.305.89"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "2.184.30"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.74%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0935"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "2.510.28"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.865"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.27%  "
$ws.Range("D18").Value = "2.193.67"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "41.169.82"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "0.0₃0948"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +20.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.98%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.122"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0725"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.56"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0297"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.201"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.79%  "
$ws.Range("E47").Value = "  +2.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +3.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.15%  "
